$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename sheets to unify the DataNode / DataTable / Entity naming convention.
$wb.Worksheets.Item("Property1").Name = "DataNode_1"
$wb.Worksheets.Item("Property2").Name = "DataNode_2"
$wb.Worksheets.Item("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets.Item("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets.Item("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets.Item("Record_Task").Name = "DataTable_Task"

# Remove the now-unused Record_Building sheet (and its comments/strings).
$null = $wb.Worksheets.Item("Record_Building").Delete()

# Make the Hero table the active sheet (matches the saved view state).
$wb.Worksheets.Item("DataTable_Hero").Activate()

$excel.DisplayAlerts = $true
